# Replace the two FESCO bill rows (rows 2 and 3) with a single GEPCO
# bill row, matching the target sharedStrings.xml / sheet1.xml produced
# by the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two existing data rows entirely.
$ws.Rows.Item(2).EntireRow.Delete()
$ws.Rows.Item(2).EntireRow.Delete()

# Helper: write a value into a cell as a text formula ( ="..." ), which
# we later convert to a plain (non-formula) value via copy / paste
# special. This keeps Excel from "helpfully" re-interpreting
# numeric-looking or date-looking strings (e.g. "30123630360500",
# "800", "AUG-2025") as numbers/dates - the source workbook stores
# every one of these as a plain shared string - without touching
# NumberFormat/Style (which would otherwise leave extra, unused style
# entries behind in styles.xml).
function Set-TextFormula($cellRef, $value) {
    $escaped = $value -replace '"', '""'
    $ws.Range($cellRef).Formula = '="' + $escaped + '"'
}

Set-TextFormula "C2" "30123630360500"
Set-TextFormula "D2" "30123630360500"
Set-TextFormula "E2" "800"
Set-TextFormula "F2" "800"
Set-TextFormula "H2" "AUG-2025"
Set-TextFormula "I2" "1552419"
Set-TextFormula "J2" "2570"
Set-TextFormula "K2" "21760"
Set-TextFormula "L2" "320"
Set-TextFormula "M2" "3480"
Set-TextFormula "N2" "222"
Set-TextFormula "O2" "222"
Set-TextFormula "Q2" "35.1500"
Set-TextFormula "R2" "43.8200"

# Convert each contiguous block of formulas into static text values.
$block1 = $ws.Range("C2:F2")
$block1.Copy()
$block1.PasteSpecial(-4163)   # xlPasteValues

$block2 = $ws.Range("H2:O2")
$block2.Copy()
$block2.PasteSpecial(-4163)

$block3 = $ws.Range("Q2:R2")
$block3.Copy()
$block3.PasteSpecial(-4163)

# These values are never mistaken for numbers/dates, so they can be
# assigned directly.
$ws.Range("A2").Value = "GEPCO_30123630360500.pdf"
$ws.Range("B2").Value = "GEPCO"
$ws.Range("G2").Value = "A-2c(06)T"
$ws.Range("S2").Value = "D:\BILLs App\BillWebApp\BillWebApp\downloads\GEPCO_30123630360500.pdf"
